$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for years 2010-2025 (replaces old 2000-2025 data)
$data = @(
    @(2010, 8.8000000000000007, 6.3),
    @(2011, 6.1, 8.8000000000000007),
    @(2012, 6.6, 10.65),
    @(2013, 6.5, 10.119999999999999),
    @(2014, 11.4, 8.31),
    @(2015, 12.9, 11.4),
    @(2016, 5.4, 4),
    @(2017, 2.5, 5.8),
    @(2018, 2.2999999999999998, 3.7),
    @(2019, 3, 7.05),
    @(2020, 4.9000000000000004, 6.6),
    @(2021, 8.4, 6.3),
    @(2022, 11.9, 19.46),
    @(2023, 7.4, 4.8),
    @(2024, 9.51, 7.5),
    @(2025, 5.26, 9.5)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# The "highlighted" indexation cell for 2017 used to live at C19 (old layout);
# with the table now starting at 2010, 2017's row is C9. Bring its
# distinct formatting along before the old row is removed.
$ws.Range("C19").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Remove the now-unused trailing rows (previously rows 18-27), shifting cells up
$ws.Range("A18:C27").Delete(-4162)

# Keep the sorted-range bookkeeping in sync with the now-smaller table
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A1"))
$ws.Sort.SetRange($ws.Range("A2:C18"))
$ws.Sort.Header = -4105
$ws.Sort.Apply()

# Update selection to match target
$ws.Range("C5").Select()
